# Apply updated dSF (column F) values as re-pulled from source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new dSF (column F) value.
$updates = @{
    3  = -3
    6  = -2
    7  = 1
    12 = -5
    13 = 9
    14 = 2
    16 = 2
    20 = 2
    22 = 4
    23 = -5
    25 = 1
    28 = 6
    31 = 2
    34 = 1
    37 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
